# Auto-generated Excel COM-interop edit script
# Applies odds updates to rows 3, 6, 8, 9 and appends new row 10 (Al Ittihad vs Al Ahli SC)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 odds updates ---
$ws.Range("G3").Value = 4.15
$ws.Range("I3").Value = 2.18
$ws.Range("J3").Value = 4.75
$ws.Range("K3").Value = 1.83
$ws.Range("L3").Value = 2.82
$ws.Range("Q3").Value = 2.8
$ws.Range("S3").Value = 1.6
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 2.22
$ws.Range("W3").Value = 7.8
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 80
$ws.Range("AA3").Value = 55
$ws.Range("AB3").Value = 75
$ws.Range("AE3").Value = 19
$ws.Range("AH3").Value = 5.2
$ws.Range("AI3").Value = 9
$ws.Range("AJ3").Value = 9.5
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 23
$ws.Range("AM3").Value = 45
$ws.Range("AN3").Value = 5.7
$ws.Range("AO3").Value = 27
$ws.Range("AP3").Value = 37
$ws.Range("AQ3").Value = 175
$ws.Range("AR3").Value = 250
$ws.Range("AW3").Value = 3.8
$ws.Range("AX3").Value = 11.75
$ws.Range("AY3").Value = 24
$ws.Range("AZ3").Value = 55
$ws.Range("BA3").Value = 100

# --- Row 6 odds updates ---
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3.65
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 2.45
$ws.Range("K6").Value = 2.4
$ws.Range("L6").Value = 3.5
$ws.Range("N6").Value = 9.5
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.8
$ws.Range("Q6").Value = 1.45
$ws.Range("R6").Value = 2.55
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.6
$ws.Range("U6").Value = 1.42
$ws.Range("V6").Value = 2.67
$ws.Range("W6").Value = 12.5
$ws.Range("X6").Value = 13.5
$ws.Range("Y6").Value = 8.75
$ws.Range("Z6").Value = 21
$ws.Range("AA6").Value = 13.5
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 7.9
$ws.Range("AE6").Value = 10.75
$ws.Range("AF6").Value = 30
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 11.5
$ws.Range("AK6").Value = 45
$ws.Range("AL6").Value = 24
$ws.Range("AM6").Value = 22
$ws.Range("AN6").Value = 4.55
$ws.Range("AO6").Value = 9.75
$ws.Range("AP6").Value = 13
$ws.Range("AQ6").Value = 30
$ws.Range("AR6").Value = 40
$ws.Range("AS6").Value = 100
$ws.Range("AT6").Value = 3.6
$ws.Range("AW6").Value = 5.9
$ws.Range("AX6").Value = 16.5
$ws.Range("AY6").Value = 17
$ws.Range("AZ6").Value = 65
$ws.Range("BA6").Value = 70
$ws.Range("BB6").Value = 150
$ws.Range("BC6").Value = 400

# --- Row 8 odds updates ---
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 5.25
$ws.Range("I8").Value = 6.5
$ws.Range("J8").Value = 1.83
$ws.Range("K8").Value = 2.5
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 15
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.33
$ws.Range("Q8").Value = 1.6
$ws.Range("R8").Value = 2.3
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.4
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("W8").Value = 8
$ws.Range("Z8").Value = 9
$ws.Range("AC8").Value = 15
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 800
$ws.Range("AH8").Value = 19
$ws.Range("AJ8").Value = 21
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 6.5
$ws.Range("AQ8").Value = 17
$ws.Range("AS8").Value = 101
$ws.Range("AT8").Value = 3.4
$ws.Range("BA8").Value = 126

# --- Row 9 odds updates ---
$ws.Range("G9").Value = 1.8
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 3.9
$ws.Range("J9").Value = 2.38
$ws.Range("K9").Value = 2.25
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 8
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 12
$ws.Range("AD9").Value = 7.5
$ws.Range("AO9").Value = 9.5
$ws.Range("AP9").Value = 19
$ws.Range("AQ9").Value = 29

# --- New row 10: Al Ittihad vs Al Ahli SC ---
$ws.Range("A10").Value = "G2TuVbho"
$ws.Range("B10").Value = "31/10/2024"
$ws.Range("C10").Value = "15:00"
$ws.Range("D10").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E10").Value = "Al Ittihad"
$ws.Range("F10").Value = "Al Ahli SC"
$ws.Range("G10").Value = 2.15
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 2.88
$ws.Range("J10").Value = 2.6
$ws.Range("K10").Value = 2.5
$ws.Range("L10").Value = 3.1
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 11
$ws.Range("O10").Value = 1.13
$ws.Range("P10").Value = 5.5
$ws.Range("Q10").Value = 1.44
$ws.Range("R10").Value = 2.63
$ws.Range("S10").Value = 1.25
$ws.Range("T10").Value = 3.75
$ws.Range("U10").Value = 1.4
$ws.Range("V10").Value = 2.75
$ws.Range("W10").Value = 13
$ws.Range("X10").Value = 15
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 15
$ws.Range("AB10").Value = 19
$ws.Range("AC10").Value = 21
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 11
$ws.Range("AF10").Value = 29
$ws.Range("AG10").Value = 81
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 12
$ws.Range("AK10").Value = 34
$ws.Range("AL10").Value = 21
$ws.Range("AM10").Value = 21
$ws.Range("AN10").Value = 4.75
$ws.Range("AO10").Value = 11
$ws.Range("AP10").Value = 17
$ws.Range("AQ10").Value = 34
$ws.Range("AR10").Value = 41
$ws.Range("AS10").Value = 81
$ws.Range("AT10").Value = 3.75
$ws.Range("AU10").Value = 7
$ws.Range("AV10").Value = 41
$ws.Range("AW10").Value = 5.5
$ws.Range("AX10").Value = 15
$ws.Range("AY10").Value = 19
$ws.Range("AZ10").Value = 41
$ws.Range("BA10").Value = 51
$ws.Range("BB10").Value = 81
$ws.Range("BC10").Value = 300
$ws.Range("BD10").Value = 81
